$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 405.1
$ws.Range("J2").Value = 395.66666
$ws.Range("L2").Value = 395.66666
$ws.Range("N2").Value = -621.66666
# Row 5
$ws.Range("H5").Value = 243.72728
$ws.Range("I5").Value = 118.6
$ws.Range("J5").Value = 348
$ws.Range("K5").Value = 118.6
$ws.Range("L5").Value = 348
$ws.Range("M5").Value = -3.599999999999994
$ws.Range("N5").Value = -578
# Row 33
$ws.Range("H33").Value = 225.55556
$ws.Range("I33").Value = 154.09091
$ws.Range("J33").Value = 337.85715
$ws.Range("K33").Value = 154.09091
$ws.Range("L33").Value = 337.85715
$ws.Range("M33").Value = 74.90908999999999
$ws.Range("N33").Value = -795.85715
# Row 62
$ws.Range("H62").Value = 125001660
$ws.Range("J62").Value = 3333
$ws.Range("L62").Value = 3333
$ws.Range("N62").Value = -4581
# Row 65
$ws.Range("H65").Value = 125001660
$ws.Range("J65").Value = 3333
$ws.Range("L65").Value = 16665
$ws.Range("N65").Value = -22905
# Row 98
$ws.Range("H98").Value = 9608.916999999999
$ws.Range("I98").Value = 788.625
$ws.Range("K98").Value = 788.625
$ws.Range("M98").Value = 709.375
# Row 111
$ws.Range("H111").Value = 2335.4546
$ws.Range("I111").Value = 2417.5557
$ws.Range("J111").Value = 1966
$ws.Range("K111").Value = 7252.6671
$ws.Range("L111").Value = 5898
$ws.Range("M111").Value = -4185.6671
$ws.Range("N111").Value = -12032
# Row 116
$ws.Range("H116").Value = 6261.6924
$ws.Range("I116").Value = 6976.7
$ws.Range("K116").Value = 6976.7
$ws.Range("M116").Value = -3534.7
# Row 122
$ws.Range("H122").Value = 9608.916999999999
$ws.Range("I122").Value = 788.625
$ws.Range("K122").Value = 2365.875
$ws.Range("M122").Value = 84.125
# Row 127
$ws.Range("H127").Value = 711
$ws.Range("I127").Value = 711
$ws.Range("K127").Value = 2133
$ws.Range("M127").Value = 2827
# Row 129
$ws.Range("H129").Value = 1666.7693
$ws.Range("I129").Value = 1007.7778
$ws.Range("J129").Value = 3149.5
$ws.Range("K129").Value = 3023.3334
$ws.Range("L129").Value = 9448.5
$ws.Range("M129").Value = 1976.6666
$ws.Range("N129").Value = -19448.5
# Row 135
$ws.Range("H135").Value = 2058.25
$ws.Range("I135").Value = 1799.9546
$ws.Range("K135").Value = 16199.5914
$ws.Range("M135").Value = -13664.5914
# Row 138
$ws.Range("H138").Value = 3321.8545
$ws.Range("I138").Value = 4477.5293
$ws.Range("J138").Value = 2804.842
$ws.Range("K138").Value = 13432.5879
$ws.Range("L138").Value = 8414.526
$ws.Range("M138").Value = -8292.5879
$ws.Range("N138").Value = -18694.526

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 9141.076999999999
$ws.Range("I88").Value = 1389.8334
$ws.Range("J88").Value = 15785
$ws.Range("K88").Value = 1389.8334
$ws.Range("L88").Value = 15785
$ws.Range("M88").Value = -983.8334
$ws.Range("N88").Value = -16597
# Row 91
$ws.Range("H91").Value = 9141.076999999999
$ws.Range("I91").Value = 1389.8334
$ws.Range("J91").Value = 15785
$ws.Range("K91").Value = 1389.8334
$ws.Range("L91").Value = 15785
$ws.Range("M91").Value = 14.16660000000002
$ws.Range("N91").Value = -18593
# Row 122
$ws.Range("H122").Value = 3457.182
$ws.Range("I122").Value = 2106
$ws.Range("K122").Value = 6318
$ws.Range("M122").Value = -3868

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 1672.3334
$ws.Range("I94").Value = 1250.1428
$ws.Range("J94").Value = 2263.4
$ws.Range("K94").Value = 1250.1428
$ws.Range("L94").Value = 2263.4
$ws.Range("M94").Value = -799.1428000000001
$ws.Range("N94").Value = -3165.4
# Row 99
$ws.Range("H99").Value = 1848.0667
$ws.Range("I99").Value = 1760.3334
$ws.Range("K99").Value = 1760.3334
$ws.Range("M99").Value = -262.3334
# Row 105
$ws.Range("H105").Value = 1794.4286
$ws.Range("I105").Value = 1713.2
$ws.Range("K105").Value = 1713.2
$ws.Range("M105").Value = 33.79999999999995
# Row 126
$ws.Range("H126").Value = 1848.0667
$ws.Range("I126").Value = 1760.3334
$ws.Range("K126").Value = 5281.0002
$ws.Range("M126").Value = -2811.0002
# Row 141
$ws.Range("H141").Value = 158312.9
$ws.Range("J141").Value = 178115.47
$ws.Range("L141").Value = 178115.47
$ws.Range("N141").Value = -188475.47

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 113
$ws.Range("H113").Value = 1008.7273
$ws.Range("I113").Value = 1305
$ws.Range("J113").Value = 653.2
$ws.Range("K113").Value = 3915
$ws.Range("L113").Value = 1959.6
$ws.Range("M113").Value = -1745
$ws.Range("N113").Value = -6299.6
# Row 120
$ws.Range("H120").Value = 31446.076
$ws.Range("J120").Value = 38888.777
$ws.Range("L120").Value = 116666.331
$ws.Range("N120").Value = -126342.331

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 1364.6666
$ws.Range("I3").Value = 1772
$ws.Range("J3").Value = 550
$ws.Range("K3").Value = 1772
$ws.Range("L3").Value = 550
$ws.Range("M3").Value = -1656
$ws.Range("N3").Value = -782
# Row 14
$ws.Range("H14").Value = 1753202.1
$ws.Range("I14").Value = 2627052
$ws.Range("J14").Value = 5502.5
$ws.Range("K14").Value = 2627052
$ws.Range("L14").Value = 5502.5
$ws.Range("M14").Value = -2626884
$ws.Range("N14").Value = -5838.5
# Row 24
$ws.Range("H24").Value = 16423.334
$ws.Range("I24").Value = 15000
$ws.Range("J24").Value = 17846.666
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 17846.666
$ws.Range("M24").Value = -14827
$ws.Range("N24").Value = -18192.666
# Row 97
$ws.Range("H97").Value = 2748.4443
$ws.Range("I97").Value = 540.9167
$ws.Range("J97").Value = 7163.5
$ws.Range("K97").Value = 540.9167
$ws.Range("L97").Value = 7163.5
$ws.Range("M97").Value = -44.91669999999999
$ws.Range("N97").Value = -8155.5
# Row 102
$ws.Range("H102").Value = 8228.725
$ws.Range("I102").Value = 9198.870999999999
$ws.Range("K102").Value = 9198.870999999999
$ws.Range("M102").Value = -7576.870999999999
# Row 113
$ws.Range("H113").Value = 5645.8237
$ws.Range("I113").Value = 5452.25
$ws.Range("J113").Value = 5705.385
$ws.Range("K113").Value = 5452.25
$ws.Range("L113").Value = 5705.385
$ws.Range("M113").Value = -3282.25
$ws.Range("N113").Value = -10045.385
# Row 122
$ws.Range("H122").Value = 3997.138
$ws.Range("I122").Value = 3544.2307
$ws.Range("J122").Value = 4365.125
$ws.Range("K122").Value = 10632.6921
$ws.Range("L122").Value = 13095.375
$ws.Range("M122").Value = -8182.6921
$ws.Range("N122").Value = -17995.375
# Row 126
$ws.Range("H126").Value = 6955.826
$ws.Range("I126").Value = 6545.875
$ws.Range("K126").Value = 19637.625
$ws.Range("M126").Value = -17167.625
# Row 132
$ws.Range("H132").Value = 1639.0333
$ws.Range("I132").Value = 1038.5454
$ws.Range("K132").Value = 3115.6362
$ws.Range("M132").Value = -585.6361999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2512.0715
$ws.Range("J68").Value = 2700
$ws.Range("L68").Value = 2700
$ws.Range("N68").Value = -4198
# Row 71
$ws.Range("H71").Value = 2512.0715
$ws.Range("J71").Value = 2700
$ws.Range("L71").Value = 13500
$ws.Range("N71").Value = -20988
# Row 93
$ws.Range("H93").Value = 1386.8889
$ws.Range("I93").Value = 954.7143
$ws.Range("K93").Value = 954.7143
$ws.Range("M93").Value = 293.2857
# Row 100
$ws.Range("H100").Value = 2419.8572
$ws.Range("I100").Value = 1372.5
$ws.Range("J100").Value = 3816.3333
$ws.Range("K100").Value = 1372.5
$ws.Range("L100").Value = 3816.3333
$ws.Range("M100").Value = -831.5
$ws.Range("N100").Value = -4898.3333
# Row 101
$ws.Range("H101").Value = 2303
$ws.Range("J101").Value = 2303
$ws.Range("L101").Value = 2303
$ws.Range("N101").Value = -8793
# Row 122
$ws.Range("H122").Value = 5623.55
$ws.Range("I122").Value = 4939.4546
$ws.Range("K122").Value = 14818.3638
$ws.Range("M122").Value = -12368.3638

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 7249
$ws.Range("I14").Value = 9498
$ws.Range("K14").Value = 9498
$ws.Range("M14").Value = -9330
# Row 122
$ws.Range("H122").Value = 6240.1
$ws.Range("I122").Value = 3649.3333
$ws.Range("J122").Value = 10126.25
$ws.Range("K122").Value = 10947.9999
$ws.Range("L122").Value = 30378.75
$ws.Range("M122").Value = -8497.999899999999
$ws.Range("N122").Value = -35278.75
# Row 126
$ws.Range("H126").Value = 2251.4707
$ws.Range("I126").Value = 2219.7856
$ws.Range("J126").Value = 2399.3333
$ws.Range("K126").Value = 6659.3568
$ws.Range("L126").Value = 7197.999899999999
$ws.Range("M126").Value = -4189.3568
$ws.Range("N126").Value = -12137.9999
# Row 132
$ws.Range("H132").Value = 4835.7188
$ws.Range("I132").Value = 4649.76
$ws.Range("K132").Value = 13949.28
$ws.Range("M132").Value = -11419.28
